$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions): simple F-column (want-to-go count) bumps ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value2 = 144
$ws1.Range("F6").Value2 = 298
$ws1.Range("F7").Value2 = 13314
$ws1.Range("F8").Value2 = 83
$ws1.Range("F9").Value2 = 338
$ws1.Range("F10").Value2 = 5314
$ws1.Range("F21").Value2 = 3779
$ws1.Range("F22").Value2 = 115
$ws1.Range("F24").Value2 = 5014
$ws1.Range("F26").Value2 = 2004
$ws1.Range("F28").Value2 = 301
$ws1.Range("F29").Value2 = 7353
$ws1.Range("F30").Value2 = 26
$ws1.Range("F31").Value2 = 169
$ws1.Range("F32").Value2 = 2160
$ws1.Range("F33").Value2 = 2100
$ws1.Range("F34").Value2 = 1316
$ws1.Range("F35").Value2 = 135
$ws1.Range("F36").Value2 = 1136
$ws1.Range("F37").Value2 = 10
$ws1.Range("F40").Value2 = 5
$ws1.Range("F41").Value2 = 1163
$ws1.Range("F42").Value2 = 1161
$ws1.Range("F45").Value2 = 1275
$ws1.Range("F46").Value2 = 1933
$ws1.Range("F47").Value2 = 91
$ws1.Range("F48").Value2 = 182

# ---- Sheet "演出" (Performances): two listings removed, remaining rows shift up;
#      refreshed want-to-go counts applied to the surviving rows ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(14).Delete()
$ws2.Range("F2").Value2 = 48
$ws2.Range("F3").Value2 = 38
$ws2.Range("F5").Value2 = 140
for ($r = 2; $r -le 20; $r++) {
    $ws2.Cells.Item($r, 1).Value2 = $r - 1
}

# ---- Sheet "本地生活" (Local life): simple F-column bumps ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value2 = 514
$ws3.Range("F3").Value2 = 687
$ws3.Range("F4").Value2 = 54

# ---- Sheet "全部类型" (All types): simple F-column bumps ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value2 = 144
$ws4.Range("F4").Value2 = 41
$ws4.Range("F5").Value2 = 38
$ws4.Range("F6").Value2 = 514
$ws4.Range("F7").Value2 = 687
$ws4.Range("F8").Value2 = 298
$ws4.Range("F9").Value2 = 13314
$ws4.Range("F10").Value2 = 338
$ws4.Range("F11").Value2 = 5314
$ws4.Range("F18").Value2 = 140
$ws4.Range("F20").Value2 = 3779
$ws4.Range("F22").Value2 = 115
$ws4.Range("F23").Value2 = 5014
$ws4.Range("F25").Value2 = 2004
$ws4.Range("F27").Value2 = 301
$ws4.Range("F28").Value2 = 7353
$ws4.Range("F29").Value2 = 26
$ws4.Range("F30").Value2 = 169
$ws4.Range("F31").Value2 = 2160
$ws4.Range("F32").Value2 = 2100
$ws4.Range("F33").Value2 = 1316
$ws4.Range("F34").Value2 = 135
$ws4.Range("F35").Value2 = 1136
$ws4.Range("F36").Value2 = 10
$ws4.Range("F39").Value2 = 5
$ws4.Range("F40").Value2 = 1163
$ws4.Range("F41").Value2 = 1161
$ws4.Range("F45").Value2 = 1275
$ws4.Range("F46").Value2 = 1933
$ws4.Range("F47").Value2 = 91
$ws4.Range("F48").Value2 = 182

